$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-strings table was reordered (words shuffled among rows) while
# the per-row counts (column B) stayed put. Update column A text per row to
# match the new association between word and its existing count.
$ws.Range("A20").Value = 'говядина'
$ws.Range("A21").Value = 'сено'
$ws.Range("A24").Value = 'выбойка'
$ws.Range("A25").Value = 'чулок'
$ws.Range("A26").Value = 'шелк'
$ws.Range("A27").Value = 'сахар'
$ws.Range("A30").Value = 'китайка'
$ws.Range("A31").Value = 'сапог'
$ws.Range("A32").Value = 'сани'
$ws.Range("A33").Value = 'коса'
$ws.Range("A35").Value = 'конь'
$ws.Range("A36").Value = 'рогожа'
$ws.Range("A38").Value = 'веревка'
$ws.Range("A39").Value = 'платок'
$ws.Range("A40").Value = 'замок'
$ws.Range("A41").Value = 'овца'
$ws.Range("A42").Value = 'горшок'
$ws.Range("A43").Value = 'гвоздь'
$ws.Range("A44").Value = 'обод'
$ws.Range("A45").Value = 'котел'
$ws.Range("A46").Value = 'нитка'
$ws.Range("A47").Value = 'скотский кожа'
$ws.Range("A48").Value = 'сосуд'
$ws.Range("A49").Value = 'гумми'
$ws.Range("A50").Value = 'роза'
$ws.Range("A51").Value = 'брусья'
$ws.Range("A52").Value = 'покроми'
$ws.Range("A53").Value = 'хомут'
$ws.Range("A54").Value = 'сковорода'
$ws.Range("A55").Value = 'бечева'
$ws.Range("A56").Value = 'дуга'
